$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5607193805
$ws.Range("C2").Value = -224.72985331
$ws.Range("D2").Value = -225.29057269
$ws.Range("E2").Value = -224.5047077181

$ws.Range("B3").Value = -0.5693031318
$ws.Range("C3").Value = -224.66900792
$ws.Range("D3").Value = -225.23831105
$ws.Range("E3").Value = -224.5047077181

$ws.Range("B4").Value = -0.5700589527
$ws.Range("C4").Value = -224.63837187
$ws.Range("D4").Value = -225.20843082
$ws.Range("E4").Value = -224.5047077181

$wb.Save()
